# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.308.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "'1.808.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'313.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.5154"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.3976"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.11%  "
$ws.Range("D9").Value = "'0.07827"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.46%  "
$ws.Range("D10").Value = "'1.112"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "'40.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").Value = "'6.328"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'20.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "'1.810.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'7.311"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "'92.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'17.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "'6.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "'28.346.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "'11.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").Value = "'2.231"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "'161.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.020.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").Value = "'2.417"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'127.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'0.1101"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").Value = "'3.663"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "'5.567"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("D35").Value = "'0.07186"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.68%  "
$ws.Range("D36").Value = "'9.131"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("D37").Value = "'0.02356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "'5.049"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").Value = "'11.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("D41").Value = "'0.6182"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Value = "'13.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").Value = "'0.5986"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.304"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.59%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.743"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "'125.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "'1.214"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("D51").Value = "'0.06832"
$ws.Range("D51").Style = "Normal"
